$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testcases")

$ws.Range('B2').Value = 'Component: Multi-Functional Tool Application'
$ws.Range('C6').Value = 'Application is installed on a desktop computer'
$ws.Range('E6').Value = '1. Copy MultiFunctionalTool_For_Desktop.zip from specified tec-share location<br>2. Extract contents to preferred location<br>3. Double-click on MultiFunctionalToolApplication'
$ws.Range('F6').Value = 'Application launches successfully with all features accessible'
$ws.Range('C7').Value = 'Application is installed on a laptop'
$ws.Range('E7').Value = '1. Copy MultiFunctionalTool_For_Laptop.zip from specified tec-share location<br>2. Extract contents to preferred location<br>3. Double-click on MultiFunctionalToolApplication'
$ws.Range('F7').Value = 'Application launches successfully with all features accessible'
$ws.Range('C8').Value = 'Application is installed and running'
$ws.Range('D8').Value = 'Verify Network Packet Capture start functionality'
$ws.Range('E8').Value = '1. Navigate to Network Packet Capture section<br>2. Click Start button'
$ws.Range('F8').Value = 'Packet capture begins successfully'
$ws.Range('C9').Value = 'Network Packet Capture is running'
$ws.Range('D9').Value = 'Verify Network Packet Capture stop functionality'
$ws.Range('E9').Value = '1. Navigate to Network Packet Capture section<br>2. Click Stop button'
$ws.Range('F9').Value = '1. Packet capture stops<br>2. .pcap file is generated<br>3. File is copied to MFP''s Shared Folder<br>4. Shared Folder opens automatically'
$ws.Range('C10').Value = 'Application is installed and running'
$ws.Range('D10').Value = 'Verify Memory Leak Check functionality'
$ws.Range('E10').Value = '1. Navigate to Memory Leak Check section<br>2. Select a protocol<br>3. Run the memory leak check'
$ws.Range('F10').Value = 'Memory Leak Comparison Table is displayed with accurate information'
$ws.Range('C11').Value = 'Application is installed and running'
$ws.Range('D11').Value = 'Verify Debug Log Collection functionality'
$ws.Range('E11').Value = '1. Navigate to Debug Log Collection section<br>2. Click Run button'
$ws.Range('F11').Value = '1. Script executes successfully<br>2. Logs are collected<br>3. Logs are copied to MFP''s Shared Folder<br>4. Shared Folder opens automatically'
$ws.Range('C12').Value = 'Debug Log Collection has been run once with empty folder result'
$ws.Range('D12').Value = 'Verify Debug Log Collection retry functionality'
$ws.Range('E12').Value = '1. Navigate to Debug Log Collection section<br>2. Click Run button again'
$ws.Range('F12').Value = '1. Script executes successfully<br>2. Logs are collected<br>3. Logs are copied to MFP''s Shared Folder<br>4. Shared Folder opens with logs visible'
$ws.Range('C13').Value = 'Application is installed and running'
$ws.Range('D13').Value = 'Verify Diagnostic Code Details for ECC'
$ws.Range('E13').Value = '1. Navigate to Diagnostic Code Details section<br>2. Select ECC option<br>3. Choose a specific diagnostic code'
$ws.Range('F13').Value = 'Relevant job-specific details for the selected ECC diagnostic code are displayed'
$ws.Range('C14').Value = 'Application is installed and running'
$ws.Range('D14').Value = 'Verify Diagnostic Code Details for Network Protocols'
$ws.Range('E14').Value = '1. Navigate to Diagnostic Code Details section<br>2. Select Network Protocols option<br>3. Choose a specific diagnostic code'
$ws.Range('F14').Value = 'Relevant job-specific details for the selected Network Protocols diagnostic code are displayed'
$ws.Range('C15').Value = 'Application is installed and running'
$ws.Range('D15').Value = 'Verify Diagnostic Code Details for High Security Mode'
$ws.Range('E15').Value = '1. Navigate to Diagnostic Code Details section<br>2. Select High Security Mode option<br>3. Choose a specific diagnostic code'
$ws.Range('F15').Value = 'Relevant job-specific details for the selected High Security Mode diagnostic code are displayed'
$ws.Range('C16').Value = 'Application is installed and running'
$ws.Range('D16').Value = 'Verify 08 Diagnostic Code Value Get functionality'
$ws.Range('E16').Value = '1. Navigate to 08 Diagnostic Code Value section<br>2. Select a diagnostic code<br>3. Click Get button'
$ws.Range('F16').Value = 'Current value of the selected 08 diagnostic code is displayed'
$ws.Range('C17').Value = 'Application is installed and running'
$ws.Range('D17').Value = 'Verify 08 Diagnostic Code Value Set functionality'
$ws.Range('E17').Value = '1. Navigate to 08 Diagnostic Code Value section<br>2. Select a diagnostic code<br>3. Enter a new value<br>4. Click Set button'
$ws.Range('F17').Value = 'The 08 diagnostic code value is updated successfully'
$ws.Range('C18').Value = 'Application is installed and running'
$ws.Range('D18').Value = 'Verify Protocol Configuration Get functionality'
$ws.Range('E18').Value = '1. Navigate to Protocol Configuration section<br>2. Select a protocol<br>3. Click Get button'
$ws.Range('F18').Value = 'Current value of the selected protocol is displayed'
$ws.Range('C19').Value = 'Application is installed and running'
$ws.Range('D19').Value = 'Verify Protocol Configuration Set functionality'
$ws.Range('E19').Value = '1. Navigate to Protocol Configuration section<br>2. Select a protocol<br>3. Enter a new value<br>4. Click Set button'
$ws.Range('F19').Value = 'Message indicating that Set protocol values operation is not yet implemented'
$ws.Range('C20').Value = 'Application is installed and running'
$ws.Range('D20').Value = 'Verify GUI responsiveness'
$ws.Range('E20').Value = '1. Launch the application<br>2. Navigate through all sections<br>3. Interact with all UI elements'
$ws.Range('F20').Value = 'UI responds promptly to all user interactions with no noticeable lag'
$ws.Range('C21').Value = 'Application is installed and running'
$ws.Range('D21').Value = 'Verify application performance during packet capture'
$ws.Range('E21').Value = '1. Start packet capture<br>2. Perform other operations in the application simultaneously<br>3. Stop packet capture'
$ws.Range('F21').Value = 'Application remains responsive during packet capture with no performance degradation'
$ws.Range('C22').Value = 'Application is installed and running'
$ws.Range('D22').Value = 'Verify error handling for invalid inputs'
$ws.Range('E22').Value = '1. Enter invalid values in input fields<br>2. Submit the invalid data'
$ws.Range('F22').Value = 'Application displays appropriate error messages without crashing'
$ws.Range('C23').Value = 'Application is installed and running'
$ws.Range('D23').Value = 'Verify compatibility with different screen resolutions'
$ws.Range('E23').Value = '1. Run application on displays with different resolutions<br>2. Check UI layout and element visibility'
$ws.Range('F23').Value = 'UI elements are properly displayed and accessible on all tested resolutions'
$ws.Range('C24').Value = 'Application is installed and running'
$ws.Range('D24').Value = 'Verify usability for first-time users'
$ws.Range('E24').Value = '1. Have a first-time user navigate through the application<br>2. Ask them to perform basic tasks without instructions'
$ws.Range('F24').Value = 'User can successfully navigate and perform basic tasks with minimal confusion'
$ws.Range('B25').Value = 'TC020'
$ws.Range('C25').Value = 'Application is installed and running'
$ws.Range('D25').Value = 'Verify time efficiency improvement'
$ws.Range('E25').Value = '1. Measure time to perform diagnostic tasks manually<br>2. Measure time to perform same tasks using the application<br>3. Calculate time savings'
$ws.Range('F25').Value = 'Application reduces testing time by approximately 80% compared to manual methods'
$ws.Range('G25').Value = ''
$ws.Range('H25').Value = ''
